$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("E2").Value = "192.168.0.24"
$ws.Range("H6").Select()
